# LZ titrations; updated CBLS aquarium script
#
# The workbook tracks CRM (certified reference material) titration
# accuracy checks, one row per run. This edit:
#   1) fixes up the formula in the last existing data row (D86) so it is
#      a plain (non-shared) formula,
#   2) appends a new data row (87) for the 2022-07-20 run, which also
#      introduces a new "Notes" value (the CRM-opened batch note, now
#      re-typed in all caps: "CRM OPENED 20220427 DMBP") and picks up a
#      one-off font tweak (10pt Lucida Console, vertically centered) on
#      the new Batch value cell, and
#   3) moves the on-screen selection down to where the next entry would
#      go.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-key D86's formula as a plain (non-shared) formula -----------------
$ws.Range("D86").Formula = "=100*(B86-C86)/C86"

# --- Append the new row of titration data ----------------------------------
$ws.Range("A87").Value = 20220720
$ws.Range("B87").Value = 2208.66421
$ws.Range("C87").Value = 2224.47
$ws.Range("D87").Formula = "=100*(B87-C87)/C87"
$ws.Range("E87").Value = 180
$ws.Range("F87").Value = "CRM OPENED 20220427 DMBP"

# One-off font styling on the new Batch value cell (10pt Lucida Console,
# black, vertically centered) -- matches the font used elsewhere for
# instrument-pasted readings.
$fmt = $ws.Range("B87")
$fmt.Font.Name = "Lucida Console"
$fmt.Font.Size = 10
$fmt.Font.Color = 0
$fmt.Font.Family = 3
$fmt.VerticalAlignment = -4108

# --- Move the visible selection / scroll position down to the new area ----
[void]$ws.Range("F95").Select()
